$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "1.001") are preserved as text instead of being parsed as numbers,
# matching the original inline-string cell type. Style is reset back to
# "Normal" afterwards so no residual number-format styling is left on the cells.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "29.951.13"
$ws.Cells.Item(2, 5).Value = "  +0.21%  "
$ws.Cells.Item(3, 4).Value = "1.878.95"
$ws.Cells.Item(3, 5).Value = "  -0.53%  "
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).Value = "0.7420"
$ws.Cells.Item(5, 5).Value = "  -3.02%  "
$ws.Cells.Item(6, 4).Value = "242.74"
$ws.Cells.Item(6, 5).Value = "  -0.04%  "
$ws.Cells.Item(7, 4).Value = "1.002"
$ws.Cells.Item(7, 5).Value = "  +0.15%  "
$ws.Cells.Item(8, 4).Value = "0.3160"
$ws.Cells.Item(8, 5).Value = "  +0.86%  "
$ws.Cells.Item(9, 4).Value = "0.07232"
$ws.Cells.Item(9, 5).Value = "  +1.17%  "
$ws.Cells.Item(10, 4).Value = "24.77"
$ws.Cells.Item(10, 5).Value = "  -3.63%  "
$ws.Cells.Item(11, 4).Value = "0.08343"
$ws.Cells.Item(11, 5).Value = "  -2.62%  "
$ws.Cells.Item(12, 2).Value = "WrappedEther"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(12, 4).Value = "1.921.14"
$ws.Cells.Item(12, 5).Value = "  +4.06%  "
$ws.Cells.Item(13, 2).Value = "Polygon"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(13, 4).Value = "0.7523"
$ws.Cells.Item(13, 5).Value = "  -1.43%  "
$ws.Cells.Item(14, 2).Value = "Polkadot"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(14, 4).Value = "5.401"
$ws.Cells.Item(14, 5).Value = "  +0.57%  "
$ws.Cells.Item(15, 4).Value = "92.30"
$ws.Cells.Item(15, 5).Value = "  -1.66%  "
$ws.Cells.Item(16, 4).Value = "6.137"
$ws.Cells.Item(16, 5).Value = "  -0.72%  "
$ws.Cells.Item(17, 4).Value = "29.974.55"
$ws.Cells.Item(17, 5).Value = "  +0.84%  "
$ws.Cells.Item(18, 4).Value = "248.39"
$ws.Cells.Item(18, 5).Value = "  +1.82%  "
$ws.Cells.Item(19, 4).Value = "13.58"
$ws.Cells.Item(19, 5).Value = "  -1.41%  "
$ws.Cells.Item(20, 4).Value = "0.000007860"
$ws.Cells.Item(20, 5).Value = "  +0.76%  "
$ws.Cells.Item(21, 4).Value = "1.002"
$ws.Cells.Item(21, 5).Value = "  +0.18%  "
$ws.Cells.Item(22, 4).Value = "2.125.83"
$ws.Cells.Item(22, 5).Value = "  +1.77%  "
$ws.Cells.Item(23, 4).Value = "8.007"
$ws.Cells.Item(23, 5).Value = "  -0.11%  "
$ws.Cells.Item(24, 4).Value = "1.002"
$ws.Cells.Item(24, 5).Value = "  +0.10%  "
$ws.Cells.Item(25, 4).Value = "0.1555"
$ws.Cells.Item(25, 5).Value = "  -3.88%  "
$ws.Cells.Item(26, 4).Value = "9.302"
$ws.Cells.Item(26, 5).Value = "  -1.10%  "
$ws.Cells.Item(27, 4).Value = "166.27"
$ws.Cells.Item(27, 5).Value = "  +2.41%  "
$ws.Cells.Item(28, 4).Value = "18.68"
$ws.Cells.Item(28, 5).Value = "  -0.43%  "
$ws.Cells.Item(29, 4).Value = "2.033"
$ws.Cells.Item(29, 5).Value = "  -0.01%  "
$ws.Cells.Item(30, 4).Value = "1.490"
$ws.Cells.Item(30, 5).Value = "  -2.81%  "
$ws.Cells.Item(31, 4).Value = "4.597"
$ws.Cells.Item(31, 5).Value = "  +2.39%  "
$ws.Cells.Item(32, 4).Value = "1.538"
$ws.Cells.Item(32, 5).Value = "  -0.09%  "
$ws.Cells.Item(33, 4).Value = "4.220"
$ws.Cells.Item(33, 5).Value = "  +2.56%  "
$ws.Cells.Item(34, 4).Value = "0.05367"
$ws.Cells.Item(34, 5).Value = "  -1.11%  "
$ws.Cells.Item(35, 4).Value = "1.246"
$ws.Cells.Item(35, 5).Value = "  +0.36%  "
$ws.Cells.Item(36, 4).Value = "0.7533"
$ws.Cells.Item(36, 5).Value = "  +1.36%  "
$ws.Cells.Item(37, 4).Value = "1.003"
$ws.Cells.Item(37, 5).Value = "  +0.39%  "
$ws.Cells.Item(38, 4).Value = "2.710"
$ws.Cells.Item(38, 5).Value = "  +0.44%  "
$ws.Cells.Item(39, 4).Value = "0.01963"
$ws.Cells.Item(39, 5).Value = "  +0.83%  "
$ws.Cells.Item(40, 4).Value = "2.760"
$ws.Cells.Item(40, 5).Value = "  -0.71%  "
$ws.Cells.Item(41, 5).Value = "  +1.69%  "
$ws.Cells.Item(42, 4).Value = "1.124.05"
$ws.Cells.Item(42, 5).Value = "  +2.00%  "
$ws.Cells.Item(43, 4).Value = "6.155"
$ws.Cells.Item(43, 5).Value = "  +1.23%  "
$ws.Cells.Item(44, 4).Value = "72.73"
$ws.Cells.Item(44, 5).Value = "  -0.29%  "
$ws.Cells.Item(45, 4).Value = "0.8628"
$ws.Cells.Item(45, 5).Value = "  +1.14%  "
$ws.Cells.Item(46, 4).Value = "104.77"
$ws.Cells.Item(46, 5).Value = "  +1.73%  "
$ws.Cells.Item(47, 4).Value = "1.002"
$ws.Cells.Item(47, 5).Value = "  +0.19%  "
$ws.Cells.Item(48, 4).Value = "1.869"
$ws.Cells.Item(48, 5).Value = "  +0.03%  "
$ws.Cells.Item(49, 4).Value = "7.621"
$ws.Cells.Item(49, 5).Value = "  -0.47%  "
$ws.Cells.Item(50, 4).Value = "9.545"
$ws.Cells.Item(50, 5).Value = "  -1.83%  "
$ws.Cells.Item(51, 4).Value = "2.032.26"
$ws.Cells.Item(51, 5).Value = "  +1.45%  "

# Restore default styling on column D (removes the temporary text-format style)
$colD.Style = "Normal"
